{"js": "const replacements = [\n  [\"2024-05-05 Sunday\", \"2024-05-06 Monday\"],\n  [\"285\u00d77=1995\", \"701\u00d77=4907\"],\n  [\"181\u00d78=1448\", \"826\u00d76=4956\"],\n  [\"484\u00d75=2420\", \"528\u00d77=3696\"],\n  [\"998\u00d78=7984\", \"909\u00d72=1818\"],\n  [\"162\u00d79=1458\", \"710\u00d74=2840\"],\n  [\"430\u00d73=1290\", \"195\u00d76=1170\"],\n  [\"920\u00d72=1840\", \"311\u00d75=1555\"],\n  [\"781\u00d79=7029\", \"140\u00d77=980\"],\n  [\"222\u00d78=1776\", \"225\u00d74=900\"],\n  [\"517\u00d72=1034\", \"955\u00d72=1910\"],\n  [\"997\u00d75=4985\", \"155\u00d74=620\"],\n  [\"215\u00d72=430\", \"626\u00d79=5634\"],\n  [\"117\u00d75=585\", \"326\u00d76=1956\"],\n  [\"891\u00d78=7128\", \"199\u00d76=1194\"],\n  [\"578\u00d78=4624\", \"236\u00d76=1416\"],\n  [\"824\u00d73=2472\", \"423\u00d74=1692\"],\n  [\"913\u00d78=7304\", \"345\u00d72=690\"],\n  [\"109\u00d77=763\", \"490\u00d74=1960\"],\n  [\"845\u00d73=2535\", \"160\u00d75=800\"],\n  [\"239\u00d75=1195\", \"349\u00d76=2094\"],\n  [\"995\u00d72=1990\", \"239\u00d79=2151\"],\n  [\"123\u00d79=1107\", \"509\u00d78=4072\"],\n  [\"426\u00d75=2130\", \"691\u00d76=4146\"],\n  [\"800\u00d76=4800\", \"465\u00d78=3720\"],\n  [\"400\u00d72=800\", \"419\u00d75=2095\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-05-05 Sunday', '2024-05-06 Monday'),\n    @('285\u00d77=1995', '701\u00d77=4907'),\n    @('181\u00d78=1448', '826\u00d76=4956'),\n    @('484\u00d75=2420', '528\u00d77=3696'),\n    @('998\u00d78=7984', '909\u00d72=1818'),\n    @('162\u00d79=1458', '710\u00d74=2840'),\n    @('430\u00d73=1290', '195\u00d76=1170'),\n    @('920\u00d72=1840', '311\u00d75=1555'),\n    @('781\u00d79=7029', '140\u00d77=980'),\n    @('222\u00d78=1776', '225\u00d74=900'),\n    @('517\u00d72=1034', '955\u00d72=1910'),\n    @('997\u00d75=4985', '155\u00d74=620'),\n    @('215\u00d72=430', '626\u00d79=5634'),\n    @('117\u00d75=585', '326\u00d76=1956'),\n    @('891\u00d78=7128', '199\u00d76=1194'),\n    @('578\u00d78=4624', '236\u00d76=1416'),\n    @('824\u00d73=2472', '423\u00d74=1692'),\n    @('913\u00d78=7304', '345\u00d72=690'),\n    @('109\u00d77=763', '490\u00d74=1960'),\n    @('845\u00d73=2535', '160\u00d75=800'),\n    @('239\u00d75=1195', '349\u00d76=2094'),\n    @('995\u00d72=1990', '239\u00d79=2151'),\n    @('123\u00d79=1107', '509\u00d78=4072'),\n    @('426\u00d75=2130', '691\u00d76=4146'),\n    @('800\u00d76=4800', '465\u00d78=3720'),\n    @('400\u00d72=800', '419\u00d75=2095'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $result = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $result) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n\nWrite-Output \"done\""}
